$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 54.69462833333333
$ws.Cells.Item(2, 8).Value = 164.083885
$ws.Cells.Item(2, 9).Value = 0.2790924419198448
$ws.Cells.Item(2, 10).Value = 0.2790924419198448
$ws.Cells.Item(2, 13).Value = 218.721583
$ws.Cells.Item(2, 14).Value = 656.164749
$ws.Cells.Item(2, 15).Value = 0.7793342808141792
$ws.Cells.Item(2, 16).Value = 0.7793342808141792
$ws.Cells.Item(2, 17).Value = 11962.89569066332
$ws.Cells.Item(2, 18).Value = 107666.0612159699
$ws.Cells.Item(2, 19).Value = 0.2175063075042753
$ws.Cells.Item(2, 20).Value = 0.2175063075042753
$ws.Cells.Item(3, 7).Value = 54.69462833333333
$ws.Cells.Item(3, 8).Value = 164.083885
$ws.Cells.Item(3, 9).Value = 0.2790924419198448
$ws.Cells.Item(3, 10).Value = 0.2790924419198448
$ws.Cells.Item(3, 15).Value = 0.164954193449581
$ws.Cells.Item(3, 16).Value = 0.164954193449581
$ws.Cells.Item(3, 17).Value = 2532.071100366939
$ws.Cells.Item(3, 18).Value = 22788.63990330246
$ws.Cells.Item(3, 19).Value = 0.04603746865476203
$ws.Cells.Item(3, 20).Value = 0.04603746865476203
$ws.Cells.Item(4, 7).Value = 54.69462833333333
$ws.Cells.Item(4, 8).Value = 164.083885
$ws.Cells.Item(4, 9).Value = 0.2790924419198448
$ws.Cells.Item(4, 10).Value = 0.2790924419198448
$ws.Cells.Item(4, 13).Value = 8.515309999999999
$ws.Cells.Item(4, 14).Value = 25.54593
$ws.Cells.Item(4, 15).Value = 0.03034118948727519
$ws.Cells.Item(4, 16).Value = 0.03034118948727519
$ws.Cells.Item(4, 17).Value = 465.7417155931166
$ws.Cells.Item(4, 18).Value = 4191.67544033805
$ws.Cells.Item(4, 19).Value = 0.008467996664756356
$ws.Cells.Item(4, 20).Value = 0.008467996664756354
$ws.Cells.Item(5, 7).Value = 54.69462833333333
$ws.Cells.Item(5, 8).Value = 164.083885
$ws.Cells.Item(5, 9).Value = 0.2790924419198448
$ws.Cells.Item(5, 10).Value = 0.2790924419198448
$ws.Cells.Item(5, 13).Value = 7.120231
$ws.Cells.Item(5, 14).Value = 21.360693
$ws.Cells.Item(5, 15).Value = 0.02537033624896462
$ws.Cells.Item(5, 16).Value = 0.02537033624896462
$ws.Cells.Item(5, 17).Value = 389.4383881924784
$ws.Cells.Item(5, 18).Value = 3504.945493732305
$ws.Cells.Item(5, 19).Value = 0.007080669096051092
$ws.Cells.Item(5, 20).Value = 0.007080669096051091
$ws.Cells.Item(6, 7).Value = 19.32115333333334
$ws.Cells.Item(6, 8).Value = 57.96346000000001
$ws.Cells.Item(6, 9).Value = 0.09859081282432611
$ws.Cells.Item(6, 10).Value = 0.09859081282432611
$ws.Cells.Item(6, 13).Value = 218.721583
$ws.Cells.Item(6, 14).Value = 656.164749
$ws.Cells.Item(6, 15).Value = 0.7793342808141792
$ws.Cells.Item(6, 16).Value = 0.7793342808141792
$ws.Cells.Item(6, 17).Value = 4225.953242452395
$ws.Cells.Item(6, 18).Value = 38033.57918207155
$ws.Cells.Item(6, 19).Value = 0.07683520020733155
$ws.Cells.Item(6, 20).Value = 0.07683520020733155
$ws.Cells.Item(7, 7).Value = 19.32115333333334
$ws.Cells.Item(7, 8).Value = 57.96346000000001
$ws.Cells.Item(7, 9).Value = 0.09859081282432611
$ws.Cells.Item(7, 10).Value = 0.09859081282432611
$ws.Cells.Item(7, 15).Value = 0.164954193449581
$ws.Cells.Item(7, 16).Value = 0.164954193449581
$ws.Cells.Item(7, 17).Value = 894.4668877341313
$ws.Cells.Item(7, 18).Value = 8050.201989607182
$ws.Cells.Item(7, 19).Value = 0.01626296801097532
$ws.Cells.Item(7, 20).Value = 0.01626296801097532
$ws.Cells.Item(8, 7).Value = 19.32115333333334
$ws.Cells.Item(8, 8).Value = 57.96346000000001
$ws.Cells.Item(8, 9).Value = 0.09859081282432611
$ws.Cells.Item(8, 10).Value = 0.09859081282432611
$ws.Cells.Item(8, 13).Value = 8.515309999999999
$ws.Cells.Item(8, 14).Value = 25.54593
$ws.Cells.Item(8, 15).Value = 0.03034118948727519
$ws.Cells.Item(8, 16).Value = 0.03034118948727519
$ws.Cells.Item(8, 17).Value = 164.5256101908667
$ws.Cells.Item(8, 18).Value = 1480.7304917178
$ws.Cells.Item(8, 19).Value = 0.002991362533607359
$ws.Cells.Item(8, 20).Value = 0.002991362533607359
$ws.Cells.Item(9, 7).Value = 19.32115333333334
$ws.Cells.Item(9, 8).Value = 57.96346000000001
$ws.Cells.Item(9, 9).Value = 0.09859081282432611
$ws.Cells.Item(9, 10).Value = 0.09859081282432611
$ws.Cells.Item(9, 13).Value = 7.120231
$ws.Cells.Item(9, 14).Value = 21.360693
$ws.Cells.Item(9, 15).Value = 0.02537033624896462
$ws.Cells.Item(9, 16).Value = 0.02537033624896462
$ws.Cells.Item(9, 17).Value = 137.5710749197534
$ws.Cells.Item(9, 18).Value = 1238.13967427778
$ws.Cells.Item(9, 19).Value = 0.002501282072411887
$ws.Cells.Item(9, 20).Value = 0.002501282072411886
$ws.Cells.Item(10, 7).Value = 11.023718
$ws.Cells.Item(10, 8).Value = 33.071154
$ws.Cells.Item(10, 9).Value = 0.05625116157486912
$ws.Cells.Item(10, 10).Value = 0.05625116157486911
$ws.Cells.Item(10, 13).Value = 218.721583
$ws.Cells.Item(10, 14).Value = 656.164749
$ws.Cells.Item(10, 15).Value = 0.7793342808141792
$ws.Cells.Item(10, 16).Value = 0.7793342808141792
$ws.Cells.Item(10, 17).Value = 2411.125051505594
$ws.Cells.Item(10, 18).Value = 21700.12546355035
$ws.Cells.Item(10, 19).Value = 0.04383845855091282
$ws.Cells.Item(10, 20).Value = 0.04383845855091281
$ws.Cells.Item(11, 7).Value = 11.023718
$ws.Cells.Item(11, 8).Value = 33.071154
$ws.Cells.Item(11, 9).Value = 0.05625116157486912
$ws.Cells.Item(11, 10).Value = 0.05625116157486911
$ws.Cells.Item(11, 15).Value = 0.164954193449581
$ws.Cells.Item(11, 16).Value = 0.164954193449581
$ws.Cells.Item(11, 17).Value = 510.3396552268647
$ws.Cells.Item(11, 18).Value = 4593.056897041783
$ws.Cells.Item(11, 19).Value = 0.009278864988184598
$ws.Cells.Item(11, 20).Value = 0.009278864988184598
$ws.Cells.Item(12, 7).Value = 11.023718
$ws.Cells.Item(12, 8).Value = 33.071154
$ws.Cells.Item(12, 9).Value = 0.05625116157486912
$ws.Cells.Item(12, 10).Value = 0.05625116157486911
$ws.Cells.Item(12, 13).Value = 8.515309999999999
$ws.Cells.Item(12, 14).Value = 25.54593
$ws.Cells.Item(12, 15).Value = 0.03034118948727519
$ws.Cells.Item(12, 16).Value = 0.03034118948727519
$ws.Cells.Item(12, 17).Value = 93.87037612258
$ws.Cells.Item(12, 18).Value = 844.8333851032199
$ws.Cells.Item(12, 19).Value = 0.001706727152222437
$ws.Cells.Item(12, 20).Value = 0.001706727152222437
$ws.Cells.Item(13, 7).Value = 11.023718
$ws.Cells.Item(13, 8).Value = 33.071154
$ws.Cells.Item(13, 9).Value = 0.05625116157486912
$ws.Cells.Item(13, 10).Value = 0.05625116157486911
$ws.Cells.Item(13, 13).Value = 7.120231
$ws.Cells.Item(13, 14).Value = 21.360693
$ws.Cells.Item(13, 15).Value = 0.02537033624896462
$ws.Cells.Item(13, 16).Value = 0.02537033624896462
$ws.Cells.Item(13, 17).Value = 78.49141863885801
$ws.Cells.Item(13, 18).Value = 706.4227677497221
$ws.Cells.Item(13, 19).Value = 0.001427110883549268
$ws.Cells.Item(13, 20).Value = 0.001427110883549267
$ws.Cells.Item(14, 7).Value = 110.9336623333333
$ws.Cells.Item(14, 8).Value = 332.800987
$ws.Cells.Item(14, 9).Value = 0.5660655836809599
$ws.Cells.Item(14, 10).Value = 0.5660655836809599
$ws.Cells.Item(14, 13).Value = 218.721583
$ws.Cells.Item(14, 14).Value = 656.164749
$ws.Cells.Item(14, 15).Value = 0.7793342808141792
$ws.Cells.Item(14, 16).Value = 0.7793342808141792
$ws.Cells.Item(14, 17).Value = 24263.58623353414
$ws.Cells.Item(14, 18).Value = 218372.2761018072
$ws.Cells.Item(14, 19).Value = 0.4411543145516595
$ws.Cells.Item(14, 20).Value = 0.4411543145516595
$ws.Cells.Item(15, 7).Value = 110.9336623333333
$ws.Cells.Item(15, 8).Value = 332.800987
$ws.Cells.Item(15, 9).Value = 0.5660655836809599
$ws.Cells.Item(15, 10).Value = 0.5660655836809599
$ws.Cells.Item(15, 15).Value = 0.164954193449581
$ws.Cells.Item(15, 16).Value = 0.164954193449581
$ws.Cells.Item(15, 17).Value = 5135.639988998879
$ws.Cells.Item(15, 18).Value = 46220.75990098992
$ws.Cells.Item(15, 19).Value = 0.09337489179565905
$ws.Cells.Item(15, 20).Value = 0.09337489179565905
$ws.Cells.Item(16, 7).Value = 110.9336623333333
$ws.Cells.Item(16, 8).Value = 332.800987
$ws.Cells.Item(16, 9).Value = 0.5660655836809599
$ws.Cells.Item(16, 10).Value = 0.5660655836809599
$ws.Cells.Item(16, 13).Value = 8.515309999999999
$ws.Cells.Item(16, 14).Value = 25.54593
$ws.Cells.Item(16, 15).Value = 0.03034118948727519
$ws.Cells.Item(16, 16).Value = 0.03034118948727519
$ws.Cells.Item(16, 17).Value = 944.6345242036565
$ws.Cells.Item(16, 18).Value = 8501.710717832908
$ws.Cells.Item(16, 19).Value = 0.01717510313668904
$ws.Cells.Item(16, 20).Value = 0.01717510313668904
$ws.Cells.Item(17, 7).Value = 110.9336623333333
$ws.Cells.Item(17, 8).Value = 332.800987
$ws.Cells.Item(17, 9).Value = 0.5660655836809599
$ws.Cells.Item(17, 10).Value = 0.5660655836809599
$ws.Cells.Item(17, 13).Value = 7.120231
$ws.Cells.Item(17, 14).Value = 21.360693
$ws.Cells.Item(17, 15).Value = 0.02537033624896462
$ws.Cells.Item(17, 16).Value = 0.02537033624896462
$ws.Cells.Item(17, 17).Value = 789.8733014893322
$ws.Cells.Item(17, 18).Value = 7108.859713403991
$ws.Cells.Item(17, 19).Value = 0.01436127419695237
$ws.Cells.Item(17, 20).Value = 0.01436127419695237

Write-Output "Applied 184 cell updates"